$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new data row at row 362, pushing every existing record
# (old rows 362..480) down by one position (to 363..481).
$ws.Rows.Item(362).Insert()

# Populate the newly-inserted row 362 with the new record.
$ws.Cells.Item(362, 1).Value = 6
$ws.Cells.Item(362, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(362, 3).Value = "Metropolitana"
$ws.Cells.Item(362, 4).Value = 45229
$ws.Cells.Item(362, 5).Value = 13
$ws.Cells.Item(362, 6).Value = 100112026
$ws.Cells.Item(362, 7).Value = "Haba"
$ws.Cells.Item(362, 8).Value = "Sin especificar"
$ws.Cells.Item(362, 9).Value = "Primera"
$ws.Cells.Item(362, 10).Value = 910
$ws.Cells.Item(362, 11).Value = 7000
$ws.Cells.Item(362, 12).Value = 8000
$ws.Cells.Item(362, 13).Value = 7121
$ws.Cells.Item(362, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(362, 15).Value = "Región Metropolitana"
$ws.Cells.Item(362, 16).Value = 285
$ws.Cells.Item(362, 17).Value = 25
$ws.Cells.Item(362, 18).Value = "Hortaliza"
